$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.355.74"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.55"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.84%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.38"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4810"
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4061"
$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08203"
$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.019"
$ws.Range("E10").Value = "  +2.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.48"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.892.04"
$ws.Range("E12").Value = "  +1.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.031"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.204"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.01"
$ws.Range("E15").Value = "  +2.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06794"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.386.49"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.622"
$ws.Range("E22").Value = "  +2.01%  "

$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.183"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.148.76"
$ws.Range("E25").Value = "  +3.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.565"
$ws.Range("E26").Value = "  +11.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.21"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.02"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.16"
$ws.Range("E30").Value = "  +2.02%  "

$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09543"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.583"
$ws.Range("E33").Value = "  +4.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.546"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.363"
$ws.Range("E35").Value = "  -0.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02280"
$ws.Range("E36").Value = "  +1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06106"
$ws.Range("E37").Value = "  +1.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.177"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.043"
$ws.Range("E39").Value = "  +2.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5959"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.80"
$ws.Range("E41").Value = "  +7.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.278"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.387"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07605"
$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.42"
$ws.Range("E46").Value = "  +1.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5570"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.949"
$ws.Range("E48").Value = "  +2.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.52"
$ws.Range("E49").Value = "  +4.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.421"
$ws.Range("E50").Value = "  +4.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.10"
$ws.Range("E51").Value = "  +1.07%  "

